$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The date-like text "2025-10-17" would otherwise be auto-converted to a
# date serial number by the smart-typing in Range.Value; force the cell to
# Text format first so it round-trips as a plain string, then drop the
# number-format override again so no stray style sticks to the cell.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-10-17"
$ws.Range("B6").Value = "ZZZ"
$ws.Range("C6").Value = "456CDX009"
$ws.Range("D6").Value = "Anna Nagar"
$ws.Range("A6").ClearFormats()
